# Applies a row-content permutation to rows 114-127 on the active sheet.
# Each destination row ends up with the full row content (columns A:AY)
# that a particular source row used to have (row 115 is unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 114
$lastRow  = 127
$firstCol = 1   # A
$lastCol  = 51  # AY

# Mapping of destination row -> source row (content that should end up there).
$mapping = @{
    114 = 127
    115 = 115
    116 = 119
    117 = 123
    118 = 125
    119 = 120
    120 = 122
    121 = 116
    122 = 114
    123 = 117
    124 = 126
    125 = 118
    126 = 121
    127 = 124
}

# Snapshot the current ("before") values of the whole block in one shot.
$srcRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$snapshot = $srcRange.Value2

# Wipe the block; we'll write back only the cells that should actually hold data,
# so any cell that had no prior value stays empty (matching the source row).
$srcRange.ClearContents()

# Columns that hold date-like text (e.g. "2023-08-24") which Excel would
# otherwise auto-convert into a real date serial number. Force them to stay
# text, then restore the default "Normal" style so no stray formatting is
# left behind.
$textColumns = @(25, 27)   # Y, AA

for ($destRow = $firstRow; $destRow -le $lastRow; $destRow++) {
    $sourceRow = $mapping[$destRow]
    $snapRowIndex = $sourceRow - $firstRow + 1

    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $value = $snapshot[$snapRowIndex, $col]
        if ($null -ne $value) {
            $cell = $ws.Cells.Item($destRow, $col)
            if (($textColumns -contains $col) -and ($value -is [string])) {
                $cell.NumberFormat = "@"
                $cell.Value2 = $value
                $cell.Style = "Normal"
            } else {
                $cell.Value2 = $value
            }
        }
    }
}
